# Adding new Test case for Notification OPQA-215
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fill TCID column first for both new rows
$ws.Cells.Item(8, 1).Value = "TestCase_E7"
$ws.Cells.Item(9, 1).Value = "TestCase_E8"

# Fill Description column for both new rows
$ws.Cells.Item(8, 3).Value = "Verify that user is able to unwatch a Patent from ALL content search results page"
$ws.Cells.Item(9, 3).Value = "Verify that user is able to unwatch a Post from ALL content search results page"

# Apply the wrap/border description-cell style (same as used on "Test Case Steps" sheet)
$ws2 = $wb.Worksheets.Item("Test Case Steps")
$ws2.Range("C6").Copy()
$ws.Range("C8:C9").PasteSpecial(-4122)

# Bold the "ALL" substring in each new description, matching the style used elsewhere in the sheet
$text8 = "Verify that user is able to unwatch a Patent from ALL content search results page"
$start8 = $text8.IndexOf("ALL") + 1
$ws.Range("C8").Characters($start8, 3).Font.Bold = $true

$text9 = "Verify that user is able to unwatch a Post from ALL content search results page"
$start9 = $text9.IndexOf("ALL") + 1
$ws.Range("C9").Characters($start9, 3).Font.Bold = $true

# Fill Jira id column for both new rows
$ws.Cells.Item(8, 2).Value = "OPQA-265"
$ws.Cells.Item(9, 2).Value = "OPQA-267"

# Fill Runmode / Results columns (reuse existing values used throughout the sheet)
$ws.Cells.Item(8, 4).Value = "Y"
$ws.Cells.Item(8, 5).Value = "PASS"
$ws.Cells.Item(9, 4).Value = "Y"
$ws.Cells.Item(9, 5).Value = "PASS"

# Apply the same border/fill style used by the other body rows to the new cells
# (copy formats only, from the row directly above, so no new style records are created)
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B9").PasteSpecial(-4122)
$ws.Range("D7:E7").Copy()
$ws.Range("D8:E9").PasteSpecial(-4122)

# Match the cell/view state seen after the edit
$ws.Range("A2").Select() | Out-Null

Write-Host "Added TestCase_E7 and TestCase_E8 rows"
